# Applies the edit described in the diff for Junction_Flooding_396.xlsx:
#  - widen several data columns (mostly 7 -> 8, with a few 7->9 / 5->7 / 6->8)
#  - replace the numeric dataset in rows 2-5 with newly measured values
#  - remove the former row 6 so the sheet keeps only 4 data rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width updates ---
# Excel COM ColumnWidth is expressed in characters of the Normal style font;
# the persisted OOXML col width value is larger by the fixed 5/6-character
# padding Excel always adds, so we subtract 5/6 before assigning to land on
# the widths required by the target file.
$ws.Columns.Item(2).ColumnWidth = 7.166666666666667
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws.Columns.Item(5).ColumnWidth = 7.166666666666667
$ws.Columns.Item(6).ColumnWidth = 7.166666666666667
$ws.Columns.Item(7).ColumnWidth = 7.166666666666667
$ws.Columns.Item(8).ColumnWidth = 7.166666666666667
$ws.Columns.Item(9).ColumnWidth = 7.166666666666667
$ws.Columns.Item(11).ColumnWidth = 7.166666666666667
$ws.Columns.Item(12).ColumnWidth = 7.166666666666667
$ws.Columns.Item(13).ColumnWidth = 7.166666666666667
$ws.Columns.Item(15).ColumnWidth = 7.166666666666667
$ws.Columns.Item(16).ColumnWidth = 7.166666666666667
$ws.Columns.Item(17).ColumnWidth = 7.166666666666667
$ws.Columns.Item(19).ColumnWidth = 6.166666666666667
$ws.Columns.Item(20).ColumnWidth = 8.166666666666666
$ws.Columns.Item(21).ColumnWidth = 7.166666666666667
$ws.Columns.Item(22).ColumnWidth = 7.166666666666667
$ws.Columns.Item(23).ColumnWidth = 7.166666666666667
$ws.Columns.Item(24).ColumnWidth = 7.166666666666667
$ws.Columns.Item(25).ColumnWidth = 6.166666666666667
$ws.Columns.Item(26).ColumnWidth = 7.166666666666667
$ws.Columns.Item(27).ColumnWidth = 7.166666666666667
$ws.Columns.Item(28).ColumnWidth = 7.166666666666667
$ws.Columns.Item(29).ColumnWidth = 7.166666666666667
$ws.Columns.Item(30).ColumnWidth = 7.166666666666667
$ws.Columns.Item(32).ColumnWidth = 7.166666666666667
$ws.Columns.Item(34).ColumnWidth = 7.166666666666667

# --- Replace data rows 2-5 with the new dataset ---
# Row 2
$ws.Cells.Item(2,1).Value = 45077.50694444445
$ws.Cells.Item(2,2).Value = 18.256
$ws.Cells.Item(2,3).Value = 12.153
$ws.Cells.Item(2,4).Value = 4.015
$ws.Cells.Item(2,5).Value = 38.832
$ws.Cells.Item(2,6).Value = 30.967
$ws.Cells.Item(2,7).Value = 14.367
$ws.Cells.Item(2,8).Value = 45.024
$ws.Cells.Item(2,9).Value = 22.106
$ws.Cells.Item(2,10).Value = 9.170999999999999
$ws.Cells.Item(2,11).Value = 13.751
$ws.Cells.Item(2,12).Value = 15.281
$ws.Cells.Item(2,13).Value = 15.885
$ws.Cells.Item(2,14).Value = 4.586
$ws.Cells.Item(2,15).Value = 14.287
$ws.Cells.Item(2,16).Value = 19.896
$ws.Cells.Item(2,17).Value = 12.285
$ws.Cells.Item(2,18).Value = 3.423
$ws.Cells.Item(2,19).Value = 2.238
$ws.Cells.Item(2,20).Value = 210.098
$ws.Cells.Item(2,21).Value = 39.697
$ws.Cells.Item(2,22).Value = 13.187
$ws.Cells.Item(2,23).Value = 26.04
$ws.Cells.Item(2,24).Value = 13.224
$ws.Cells.Item(2,25).Value = 3.103
$ws.Cells.Item(2,26).Value = 22.935
$ws.Cells.Item(2,27).Value = 11.648
$ws.Cells.Item(2,28).Value = 10.595
$ws.Cells.Item(2,29).Value = 12.441
$ws.Cells.Item(2,30).Value = 15.796
$ws.Cells.Item(2,31).Value = 3.454
$ws.Cells.Item(2,32).Value = 39.912
$ws.Cells.Item(2,33).Value = 7.199
$ws.Cells.Item(2,34).Value = 16.486

# Row 3
$ws.Cells.Item(3,1).Value = 45077.51388888889
$ws.Cells.Item(3,2).Value = 13.932
$ws.Cells.Item(3,3).Value = 9.747999999999999
$ws.Cells.Item(3,4).Value = 1.728
$ws.Cells.Item(3,5).Value = 30.113
$ws.Cells.Item(3,6).Value = 24.202
$ws.Cells.Item(3,7).Value = 10.965
$ws.Cells.Item(3,8).Value = 42.761
$ws.Cells.Item(3,9).Value = 16.87
$ws.Cells.Item(3,10).Value = 7.242
$ws.Cells.Item(3,11).Value = 10.624
$ws.Cells.Item(3,12).Value = 12.037
$ws.Cells.Item(3,13).Value = 12.557
$ws.Cells.Item(3,14).Value = 3.503
$ws.Cells.Item(3,15).Value = 10.903
$ws.Cells.Item(3,16).Value = 15.328
$ws.Cells.Item(3,17).Value = 9.462999999999999
$ws.Cells.Item(3,18).Value = 1.522
$ws.Cells.Item(3,19).Value = 1.002
$ws.Cells.Item(3,20).Value = 158.637
$ws.Cells.Item(3,21).Value = 30.538
$ws.Cells.Item(3,22).Value = 10.064
$ws.Cells.Item(3,23).Value = 20.168
$ws.Cells.Item(3,24).Value = 10.441
$ws.Cells.Item(3,25).Value = 2.057
$ws.Cells.Item(3,26).Value = 20.754
$ws.Cells.Item(3,27).Value = 8.888999999999999
$ws.Cells.Item(3,28).Value = 8.077999999999999
$ws.Cells.Item(3,29).Value = 9.468999999999999
$ws.Cells.Item(3,30).Value = 12.496
$ws.Cells.Item(3,31).Value = 1.265
$ws.Cells.Item(3,32).Value = 38.957
$ws.Cells.Item(3,33).Value = 5.523
$ws.Cells.Item(3,34).Value = 12.582

# Row 4
$ws.Cells.Item(4,1).Value = 45077.52083333334
$ws.Cells.Item(4,2).Value = 14.893
$ws.Cells.Item(4,3).Value = 10.723
$ws.Cells.Item(4,4).Value = 1.271
$ws.Cells.Item(4,5).Value = 32.312
$ws.Cells.Item(4,6).Value = 26.205
$ws.Cells.Item(4,7).Value = 11.72
$ws.Cells.Item(4,8).Value = 45.52
$ws.Cells.Item(4,9).Value = 18.033
$ws.Cells.Item(4,10).Value = 7.88
$ws.Cells.Item(4,11).Value = 11.587
$ws.Cells.Item(4,12).Value = 12.959
$ws.Cells.Item(4,13).Value = 13.59
$ws.Cells.Item(4,14).Value = 3.744
$ws.Cells.Item(4,15).Value = 11.655
$ws.Cells.Item(4,16).Value = 16.48
$ws.Cells.Item(4,17).Value = 9.99
$ws.Cells.Item(4,18).Value = 1.031
$ws.Cells.Item(4,19).Value = 0.779
$ws.Cells.Item(4,20).Value = 170.074
$ws.Cells.Item(4,21).Value = 32.598
$ws.Cells.Item(4,22).Value = 10.758
$ws.Cells.Item(4,23).Value = 21.713
$ws.Cells.Item(4,24).Value = 11.335
$ws.Cells.Item(4,25).Value = 1.952
$ws.Cells.Item(4,26).Value = 22.019
$ws.Cells.Item(4,27).Value = 9.502000000000001
$ws.Cells.Item(4,28).Value = 8.535
$ws.Cells.Item(4,29).Value = 10.016
$ws.Cells.Item(4,30).Value = 13.527
$ws.Cells.Item(4,31).Value = 0.773
$ws.Cells.Item(4,32).Value = 41.285
$ws.Cells.Item(4,33).Value = 5.972
$ws.Cells.Item(4,34).Value = 13.45

# Row 5
$ws.Cells.Item(5,1).Value = 45077.52777777778
$ws.Cells.Item(5,2).Value = 24.02
$ws.Cells.Item(5,3).Value = 17.7
$ws.Cells.Item(5,4).Value = 1.36
$ws.Cells.Item(5,5).Value = 52.21
$ws.Cells.Item(5,6).Value = 42.76
$ws.Cells.Item(5,7).Value = 18.9
$ws.Cells.Item(5,8).Value = 71.94
$ws.Cells.Item(5,9).Value = 29.09
$ws.Cells.Item(5,10).Value = 12.91
$ws.Cells.Item(5,11).Value = 19.11
$ws.Cells.Item(5,12).Value = 20.95
$ws.Cells.Item(5,13).Value = 22.08
$ws.Cells.Item(5,14).Value = 6.04
$ws.Cells.Item(5,15).Value = 18.8
$ws.Cells.Item(5,16).Value = 26.74
$ws.Cells.Item(5,17).Value = 15.84
$ws.Cells.Item(5,18).Value = 0.87
$ws.Cells.Item(5,19).Value = 0.9399999999999999
$ws.Cells.Item(5,20).Value = 278.82
$ws.Cells.Item(5,21).Value = 52.55
$ws.Cells.Item(5,22).Value = 17.35
$ws.Cells.Item(5,23).Value = 35.32
$ws.Cells.Item(5,24).Value = 18.51
$ws.Cells.Item(5,25).Value = 2.79
$ws.Cells.Item(5,26).Value = 35.21
$ws.Cells.Item(5,27).Value = 15.33
$ws.Cells.Item(5,28).Value = 13.6
$ws.Cells.Item(5,29).Value = 15.98
$ws.Cells.Item(5,30).Value = 21.96
$ws.Cells.Item(5,31).Value = 0.5600000000000001
$ws.Cells.Item(5,32).Value = 65.27
$ws.Cells.Item(5,33).Value = 9.77
$ws.Cells.Item(5,34).Value = 21.69

# --- Remove the old row 6 (dataset now has only 4 rows) ---
$ws.Rows.Item(6).Delete()
